$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued (numeric-looking string) columns on row 20 -- use a leading
# apostrophe so the engine stores them as literal text (matching the
# original inlineStr cells) instead of inferring a Number, then reset the
# cell style back to Normal so the forced quote-prefix formatting doesn't
# linger as a spurious style change.
$ws.Range("E20").Value = "'127.19000000"
$ws.Range("E20").Style = "Normal"

$ws.Range("F20").Value = "'84488.07392000"
$ws.Range("F20").Style = "Normal"

$ws.Range("H20").Value = "'10762496.22348290"
$ws.Range("H20").Style = "Normal"

$ws.Range("J20").Value = "'43243.29299000"
$ws.Range("J20").Style = "Normal"

$ws.Range("K20").Value = "'5510042.32479330"
$ws.Range("K20").Style = "Normal"

# Numeric columns on row 20
$ws.Range("I20").Value = 37823

$ws.Range("M20").Value = 127.1899999999999
$ws.Range("N20").Value = 126.74
$ws.Range("O20").Value = 127.3985714285714
$ws.Range("P20").Value = 129.7806666666667

$ws.Range("R20").Value = 127.19
$ws.Range("S20").Value = 126.8297836628268
$ws.Range("T20").Value = 128.5035265468364
$ws.Range("U20").Value = 130.4389303001382
$ws.Range("V20").Value = -1.935403753301813
$ws.Range("W20").Value = -1.810476564548283
$ws.Range("X20").Value = -0.1249271887535299
